$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "week"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "topic and materials"
$ws.Range("D1").Value = "homework"
$ws.Range("E1").Value = "notes"
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "F 09/05"
$ws.Range("C2").Value = "[Intro to coastal monitoring](https://smr-monitoring.github.io/lesson-plans/01-introductions.html)"
$ws.Range("E2").Value = "be prepared to go outside"
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "F 09/12"
$ws.Range("C3").Value = "[Water quality criteria and SAV](https://smr-monitoring.github.io/lesson-plans/02-wq-criteria.html)"
$ws.Range("D3").Value = "HW 1"
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "F 09/19"
$ws.Range("C4").Value = "SAV monitoring - St. Inigoes (be prepared to wade in the water)"
$ws.Range("D4").Value = "HW 2"
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "F 09/26"
$ws.Range("C5").Value = "[Calibration](https://smr-monitoring.github.io/lesson-plans/03-calibration.html)"
$ws.Range("D5").Value = "HW 3"
$ws.Range("E5").Value = "bring your line to class"
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "F 10/03"
$ws.Range("C6").Value = "[Continuous monitoring - practice deployment](https://smr-monitoring.github.io/lesson-plans/05-conmon.html)"
$ws.Range("D6").Value = "HW 4"
$ws.Range("E6").Value = "bring your line to class"
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "F 10/10"
$ws.Range("C7").Value = "[Continuous monitoring - data corrections](https://smr-monitoring.github.io/lesson-plans/06-conmon-drift-corrections.html)"
$ws.Range("D7").Value = "HW 5"
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "F 10/17"
$ws.Range("C8").Value = "[CONMON deploy](https://smr-monitoring.github.io/lesson-plans/07-conmon-deploy.html)"
$ws.Range("D8").Value = "HW 6"
$ws.Range("E8").Value = "be prepared to go outside"
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "F 10/24"
$ws.Range("C9").Value = "[CONMON retrieve](https://smr-monitoring.github.io/lesson-plans/08-conmon-retrieve.html)"
$ws.Range("D9").Value = "HW 7"
$ws.Range("E9").Value = "be prepared to go outside"
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "F 10/31"
$ws.Range("C10").Value = "TBA"
$ws.Range("D10").Value = "HW 8"
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "F 11/07"
$ws.Range("C11").Value = "[Discrete monitoring](https://smr-monitoring.github.io/lesson-plans/09-discrete-monitoring.html)"
$ws.Range("D11").Value = "HW 9"
$ws.Range("E11").Value = "be prepared to go outside"
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "F 11/14"
$ws.Range("C12").Value = "[Discrete monitoring](https://smr-monitoring.github.io/lesson-plans/10-discrete-boat.html)"
$ws.Range("D12").Value = "HW 10"
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "F 11/21"
$ws.Range("C13").Value = "[Skills review](https://smr-monitoring.github.io/lesson-plans/11-skills-review.html)"
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "F 11/28"
$ws.Range("C14").Value = "No class - Thanksgiving break"
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "F 12/05"
$ws.Range("C15").Value = "[Skills assessment](https://smr-monitoring.github.io/lesson-plans/12-skills-assessment.html)"
$ws.Range("E15").Value = "be prepared to go outside"
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "F 12/12"
$ws.Range("C16").Value = "DIY monitoring technology"
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "W 12/17"
$ws.Range("C17").Value = "2:00-4:15 Final written exam"

# clear removed cell (no longer has a value in the updated schedule)
$ws.Range("E14").ClearContents()

# re-apply the explicit black font color used for "topic and materials" rows
# that already had notable/boxed styling (style index 2 in the sheet)
$ws.Range("C10").Font.Color = 0
$ws.Range("C11").Font.Color = 0
$ws.Range("C12").Font.Color = 0
$ws.Range("C13").Font.Color = 0
$ws.Range("C14").Font.Color = 0
$ws.Range("C15").Font.Color = 0
$ws.Range("C17").Font.Color = 0

# move the active selection to match the author's last edited cell
[void]$ws.Range("E16").Select()
